# Generate Report for Handoff
# Refresh the "Latest Handoff Datetime" column for the three localization
# files (rows 3-5) that were just re-handed-off, in both the "zh-cn" and
# "de-de" status sheets. Row 2 / row 6 are untouched.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D3").Value = "2016-03-09 05:14:27"
$zhcn.Range("D4").Value = "2016-03-09 05:14:27"
$zhcn.Range("D5").Value = "2016-03-09 05:14:27"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D3").Value = "2016-03-09 05:14:37"
$dede.Range("D4").Value = "2016-03-09 05:14:37"
$dede.Range("D5").Value = "2016-03-09 05:14:37"
